# Update "paises" (countries) worksheet with refreshed COVID data and
# re-ranked rows (Costa Rica overtakes Bielorrusia/Portugal/Honduras/Etiopia;
# Togo overtakes Nueva Zelanda), plus updated "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 23:09"

# --- Plain value refreshes (no re-ranking, same country stays on the row) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7668362
$ws.Range("C4").Value = 31450
$ws.Range("D4").Value = 4883027
$ws.Range("E4").Value = 2570446
$ws.Range("G4").Value = 278
$ws.Range("H4").Value = 214889

# Row 6: Brasil
$ws.Range("B6").Value = 4927235
$ws.Range("C6").Value = 11946
$ws.Range("E6").Value = 517352
$ws.Range("G6").Value = 300
$ws.Range("H6").Value = 146675

# Row 13: Sudafrica
$ws.Range("B13").Value = 682215
$ws.Range("C13").Value = 926
$ws.Range("D13").Value = 615684
$ws.Range("E13").Value = 49515
$ws.Range("G13").Value = 40
$ws.Range("H13").Value = 17016

# Row 27: Israel
$ws.Range("B27").Value = 272309
$ws.Range("C27").Value = 5534
$ws.Range("D27").Value = 204355
$ws.Range("E27").Value = 66197
$ws.Range("G27").Value = 38
$ws.Range("H27").Value = 1757

# Row 87: Costa de Marfil
$ws.Range("B87").Value = 19885
$ws.Range("C87").Value = 3
$ws.Range("D87").Value = 19490
$ws.Range("E87").Value = 275

# Row 123: Republica de Yibuti
$ws.Range("B123").Value = 5421
$ws.Range("C123").Value = 2
$ws.Range("D123").Value = 5352
$ws.Range("E123").Value = 8

# Row 158: Yemen
$ws.Range("D158").Value = 1323
$ws.Range("E158").Value = 126
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 592

# Row 166: Republica del Chad
$ws.Range("B166").Value = 1223
$ws.Range("C166").Value = 6
$ws.Range("E166").Value = 62

# --- Re-ranking: Costa Rica moves up past Bielorrusia, Portugal, Honduras,
#     Etiopia (rows 51-55), pushing each of those down by one row ---

$ws.Range("A51").Value = "Costa Rica"
$ws.Range("B51").Value = 81129
$ws.Range("C51").Value = 1947
$ws.Range("D51").Value = 49703
$ws.Range("E51").Value = 30439
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 37
$ws.Range("H51").Value = 987

$ws.Range("A52").Value = "Bielorrusia"
$ws.Range("B52").Value = 80696
$ws.Range("C52").Value = 401
$ws.Range("D52").Value = 75303
$ws.Range("E52").Value = 4531
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 5
$ws.Range("H52").Value = 862

$ws.Range("A53").Value = "Portugal"
$ws.Range("B53").Value = 79885
$ws.Range("C53").Value = 734
$ws.Range("D53").Value = 50454
$ws.Range("E53").Value = 27413
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 2018

$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 79629
$ws.Range("C54").Value = 841
$ws.Range("D54").Value = 29305
$ws.Range("E54").Value = 47902
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 23
$ws.Range("H54").Value = 2422

$ws.Range("A55").Value = "Etiopia"
$ws.Range("B55").Value = 79437
$ws.Range("C55").Value = 618
$ws.Range("D55").Value = 34016
$ws.Range("E55").Value = 44191
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 8
$ws.Range("H55").Value = 1230

# --- Re-ranking: Togo overtakes Nueva Zelanda (rows 160-161 swap) ---

$ws.Range("A160").Value = "Togo"
$ws.Range("B160").Value = 1864
$ws.Range("C160").Value = 10
$ws.Range("D160").Value = 1403
$ws.Range("E160").Value = 413
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 48

$ws.Range("A161").Value = "Nueva Zelanda"
$ws.Range("B161").Value = 1855
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 1790
$ws.Range("E161").Value = 40
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 25

Write-Host "Update complete"
